# "Generate Report for Handback"
# The localization-status workbook is refreshed once the zh-cn and de-de
# handback packages have come back from translation: the Status column
# flips from "Ready for handoff" to "Handed back: in sync with en-US",
# each language sheet gets its Latest Target File / Latest Handback File /
# Latest Handback DateTime columns filled in, and the columns that now
# hold longer text are widened to fit.

$wb = $excel.ActiveWorkbook

$sourceMdName = "fb79a795-905e-4de4-91c2-ec580ac9116b.md"
$sourceMdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e31023dcd8e63693356d5bc6bdb3ede760471fcb/e2e/fb79a795-905e-4de4-91c2-ec580ac9116b.md"
$statusText   = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "Ready for handoff" -> "Handed back: in sync with en-US" ------
# (Overview mirrors each language's Status in columns E/F; each language
# sheet keeps its own Status in column C.)
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsZhCn.Range("C2").Value = $statusText
$wsDeDe.Range("C2").Value = $statusText

# --- zh-cn: Latest Target File / Latest Handback File / DateTime ----------
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $sourceMdUrl, "", "", $sourceMdName) | Out-Null
$wsZhCn.Range("J2").Value = "fb79a795-905e-4de4-91c2-ec580ac9116b.1e260ded4b1b56e70c10e3a6cc08507ffaaa9186.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-19 00:56:26"

# --- de-de: Latest Target File / Latest Handback File / DateTime ----------
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $sourceMdUrl, "", "", $sourceMdName) | Out-Null
$wsDeDe.Range("J2").Value = "fb79a795-905e-4de4-91c2-ec580ac9116b.1e260ded4b1b56e70c10e3a6cc08507ffaaa9186.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-19 00:56:33"

# --- Widen the columns that now carry longer handback text ----------------
# ColumnWidth is in "characters"; the saved OOXML `width` lands on a
# pixel-derived grid, so nudge the requested width so the stored value
# matches the target as closely as that grid allows.
function Set-ColumnWidthForStoredValue($range, [double]$targetStoredWidth) {
    $pixels = [Math]::Round($targetStoredWidth * 6)
    $range.ColumnWidth = ($pixels - 5) / 6
}

Set-ColumnWidthForStoredValue $wsOverview.Columns("E:F") 29.9777047293527

Set-ColumnWidthForStoredValue $wsZhCn.Columns("C:C") 29.9777047293527
Set-ColumnWidthForStoredValue $wsZhCn.Columns("I:J") 40

Set-ColumnWidthForStoredValue $wsDeDe.Columns("C:C") 29.9777047293527
Set-ColumnWidthForStoredValue $wsDeDe.Columns("I:J") 40
